$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2309.2307
$ws.Range("I53").Value = 1440.8334
$ws.Range("K53").Value = 1440.8334
$ws.Range("M53").Value = -803.8334
$ws.Range("H62").Value = 6182.44
$ws.Range("I62").Value = 5346.2354
$ws.Range("K62").Value = 5346.2354
$ws.Range("M62").Value = -4722.2354
$ws.Range("H65").Value = 6182.44
$ws.Range("I65").Value = 5346.2354
$ws.Range("K65").Value = 26731.177
$ws.Range("M65").Value = -23611.177
$ws.Range("H107").Value = 1008.95
$ws.Range("I107").Value = 811.6429000000001
$ws.Range("J107").Value = 1469.3334
$ws.Range("K107").Value = 811.6429000000001
$ws.Range("L107").Value = 1469.3334
$ws.Range("M107").Value = 1108.3571
$ws.Range("N107").Value = -5309.3334
$ws.Range("H110").Value = 41852.43
$ws.Range("J110").Value = 41852.43
$ws.Range("L110").Value = 41852.43
$ws.Range("N110").Value = -50032.43
$ws.Range("H112").Value = 1518227.9
$ws.Range("I112").Value = 2430.5
$ws.Range("J112").Value = 2086652
$ws.Range("K112").Value = 7291.5
$ws.Range("L112").Value = 6259956
$ws.Range("M112").Value = -6183.5
$ws.Range("N112").Value = -6262172
$ws.Range("H137").Value = 29763.42
$ws.Range("I137").Value = 41841.24
$ws.Range("K137").Value = 125523.72
$ws.Range("M137").Value = -122973.72
$ws.Range("H138").Value = 7232.8125
$ws.Range("J138").Value = 7751.291
$ws.Range("L138").Value = 23253.873
$ws.Range("N138").Value = -33533.873
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1153856
$ws.Range("I32").Value = 1415941.1
$ws.Range("J32").Value = 39994.25
$ws.Range("K32").Value = 1415941.1
$ws.Range("L32").Value = 39994.25
$ws.Range("M32").Value = -1415654.1
$ws.Range("N32").Value = -40568.25
$ws.Range("H45").Value = 3306.0667
$ws.Range("I45").Value = 3476.2307
$ws.Range("K45").Value = 3476.2307
$ws.Range("M45").Value = -3099.2307
$ws.Range("H88").Value = 1987.25
$ws.Range("I88").Value = 1749.5
$ws.Range("J88").Value = 2225
$ws.Range("K88").Value = 1749.5
$ws.Range("L88").Value = 2225
$ws.Range("M88").Value = -1343.5
$ws.Range("N88").Value = -3037
$ws.Range("H91").Value = 1987.25
$ws.Range("I91").Value = 1749.5
$ws.Range("J91").Value = 2225
$ws.Range("K91").Value = 1749.5
$ws.Range("L91").Value = 2225
$ws.Range("M91").Value = -345.5
$ws.Range("N91").Value = -5033
$ws.Range("H106").Value = 77995
$ws.Range("J106").Value = 77995
$ws.Range("L106").Value = 77995
$ws.Range("N106").Value = -80519
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17554.268
$ws.Range("I20").Value = 27232.889
$ws.Range("K20").Value = 27232.889
$ws.Range("M20").Value = -26985.889
$ws.Range("H86").Value = 5436
$ws.Range("I86").Value = 4299.625
$ws.Range("J86").Value = 8466.333000000001
$ws.Range("K86").Value = 4299.625
$ws.Range("L86").Value = 8466.333000000001
$ws.Range("M86").Value = -3176.625
$ws.Range("N86").Value = -10712.333
$ws.Range("H89").Value = 5436
$ws.Range("I89").Value = 4299.625
$ws.Range("J89").Value = 8466.333000000001
$ws.Range("K89").Value = 21498.125
$ws.Range("L89").Value = 42331.665
$ws.Range("M89").Value = -15882.125
$ws.Range("N89").Value = -53563.665
$ws.Range("H99").Value = 1994.1041
$ws.Range("I99").Value = 1878.1111
$ws.Range("K99").Value = 1878.1111
$ws.Range("M99").Value = -380.1111000000001
$ws.Range("H105").Value = 2161.923
$ws.Range("I105").Value = 1546.238
$ws.Range("J105").Value = 4747.8
$ws.Range("K105").Value = 1546.238
$ws.Range("L105").Value = 4747.8
$ws.Range("M105").Value = 200.7619999999999
$ws.Range("N105").Value = -8241.799999999999
$ws.Range("H141").Value = 90780
$ws.Range("J141").Value = 90780
$ws.Range("L141").Value = 90780
$ws.Range("N141").Value = -101140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 205333.33
$ws.Range("I4").Value = 77142.86
$ws.Range("K4").Value = 77142.86
$ws.Range("M4").Value = -77030.86
$ws.Range("H9").Value = 20555
$ws.Range("J9").Value = 20555
$ws.Range("L9").Value = 20555
$ws.Range("N9").Value = -20891
$ws.Range("H22").Value = 10165.583
$ws.Range("I22").Value = 2197.6
$ws.Range("J22").Value = 15857
$ws.Range("K22").Value = 2197.6
$ws.Range("L22").Value = 15857
$ws.Range("M22").Value = -1847.6
$ws.Range("N22").Value = -16557
$ws.Range("H107").Value = 68683.60000000001
$ws.Range("I107").Value = 85137.914
$ws.Range("K107").Value = 85137.914
$ws.Range("M107").Value = -83217.914
$ws.Range("H132").Value = 3177.0833
$ws.Range("I132").Value = 2732.5
$ws.Range("K132").Value = 8197.5
$ws.Range("M132").Value = -5667.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 4320
$ws.Range("I119").Value = 2400
$ws.Range("K119").Value = 7200
$ws.Range("M119").Value = -2362
$ws.Range("H131").Value = 10000
$ws.Range("J131").Value = 10000
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080
$ws.Range("H132").Value = 2545.1018
$ws.Range("I132").Value = 2281.5715
$ws.Range("K132").Value = 20534.1435
$ws.Range("M132").Value = -18004.1435
$ws.Range("H134").Value = 2684
$ws.Range("I134").Value = 2684
$ws.Range("K134").Value = 8052
$ws.Range("M134").Value = -2982
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2997
$ws.Range("J80").Value = 2997.5
$ws.Range("L80").Value = 2997.5
$ws.Range("N80").Value = -4993.5
$ws.Range("H83").Value = 2997
$ws.Range("J83").Value = 2997.5
$ws.Range("L83").Value = 14987.5
$ws.Range("N83").Value = -24971.5
$ws.Range("H126").Value = 6084.2856
$ws.Range("I126").Value = 6084.2856
$ws.Range("K126").Value = 18252.8568
$ws.Range("M126").Value = -15782.8568
$ws.Range("H132").Value = 7419.8945
$ws.Range("I132").Value = 7286.8486
$ws.Range("J132").Value = 8298
$ws.Range("K132").Value = 21860.5458
$ws.Range("L132").Value = 24894
$ws.Range("M132").Value = -19330.5458
$ws.Range("N132").Value = -29954
$ws.Range("H138").Value = 127886.86
$ws.Range("J138").Value = 147963.6
$ws.Range("L138").Value = 147963.6
$ws.Range("N138").Value = -158243.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5398.1377
$ws.Range("I61").Value = 5609.115
$ws.Range("K61").Value = 5609.115
$ws.Range("M61").Value = -5407.115
$ws.Range("H100").Value = 3200
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H113").Value = 5398.1377
$ws.Range("I113").Value = 5609.115
$ws.Range("K113").Value = 5609.115
$ws.Range("M113").Value = -3439.115
$ws.Range("H121").Value = 41185
$ws.Range("J121").Value = 41185
$ws.Range("L121").Value = 41185
$ws.Range("N121").Value = -44679
$ws.Range("H132").Value = 1004049.5
$ws.Range("I132").Value = 3335933
$ws.Range("K132").Value = 10007799
$ws.Range("M132").Value = -10005269
$ws.Range("H136").Value = 7557.8
$ws.Range("I136").Value = 7142.5713
$ws.Range("J136").Value = 7921.125
$ws.Range("K136").Value = 21427.7139
$ws.Range("L136").Value = 23763.375
$ws.Range("M136").Value = -18877.7139
$ws.Range("N136").Value = -28863.375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 22580.5
$ws.Range("I45").Value = 18000
$ws.Range("J45").Value = 23234.857
$ws.Range("K45").Value = 18000
$ws.Range("L45").Value = 23234.857
$ws.Range("M45").Value = -17509
$ws.Range("N45").Value = -24216.857
$ws.Range("H100").Value = 8188.3335
$ws.Range("I100").Value = 9613.714
$ws.Range("J100").Value = 3199.5
$ws.Range("K100").Value = 19227.428
$ws.Range("L100").Value = 6399
$ws.Range("M100").Value = -18686.428
$ws.Range("N100").Value = -7481
$ws.Range("H109").Value = 109500
$ws.Range("J109").Value = 109500
$ws.Range("L109").Value = 109500
$ws.Range("N109").Value = -112274
$ws.Range("H121").Value = 71959.5
$ws.Range("J121").Value = 71959.5
$ws.Range("L121").Value = 71959.5
$ws.Range("N121").Value = -75453.5
$ws.Range("H132").Value = 23839.26
$ws.Range("I132").Value = 35680.125
$ws.Range("J132").Value = 6616.1816
$ws.Range("K132").Value = 107040.375
$ws.Range("L132").Value = 19848.5448
$ws.Range("M132").Value = -104510.375
$ws.Range("N132").Value = -24908.5448
